$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15400
$ws.Range("G3").Value = 6600
$ws.Range("G4").Value = 48000
$ws.Range("G5").Value = 90000
$ws.Range("G6").Value = 90000
$ws.Range("G8").Value = 180000
$ws.Range("G9").Value = 50000
$ws.Range("G10").Value = 90000

$wb.Save()
